$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.869.38'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.706.56'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.47%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +8.17%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '656.59'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.428'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.58%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.702.50'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.89%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.58%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.398.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '96.750.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +13.83%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.702.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +6.25%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.529'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '526.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000205'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.08%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.75%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +15.19%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.186'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '32.90'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.26%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '658.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.78%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.605'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.72%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +18.78%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.163'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.21%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.974'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.92%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +18.97%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.86%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.443'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.89%  '
$ws.Range('E51').Style = 'Normal'
